$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Test: set A27 value and font properties
$ws.Range("A27").Value = "I2C JST connector"
$ws.Range("B27").Value = "C265374"
$ws.Range("C27").Value = 1.0

$ws.Range("A28").Value = "I2C JST header"
$ws.Range("B28").Value = "C265032"
$ws.Range("C28").Value = 1.0

Write-Host "done"
